$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold/centered/bordered) from H1 to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF)
$values = @{
    2  = @(8, 8)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(6, 7)
    6  = @(9, 9)
    7  = @(8, 8)
    8  = @(9, 9)
    9  = @(10, 10)
    10 = @(8, 8)
    11 = @(7, 7)
    12 = @(8, 8)
    13 = @(7, 7)
    14 = @(1, 1)
    15 = @(8, 9)
    16 = @(6, 6)
    17 = @(5, 5)
    18 = @(3, 3)
    19 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
